$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Seed row 4 cells by copying formats from analogous existing cells so the
# same style indices (date format, yellow-fill format) get reused instead
# of new style records being fabricated.
$ws.Range("B3").Copy($ws.Range("B4"))
$ws.Range("E3").Copy($ws.Range("D4"))

# Now set the actual content for the new "Day 3" row.
$ws.Range("A4").Value = "Day 3"
$ws.Range("B4").Value2 = 45805
$ws.Range("C4").Value = "Valid Parantheses"
$ws.Range("D4").Value = "Top K Frequent Elements"
$ws.Range("E4").Value = "3Sum"
$ws.Range("F4").Value = "Stack, Hash Table, Heap, Two Pointers, Sorting"
$ws.Range("G4").Value = "S"
$ws.Range("H4").Value = "YES"

# Column F grew slightly wider to fit the new content.
$ws.Columns.Item(6).ColumnWidth = 38.5

# Selection ends on H5 after the edit.
$ws.Range("H5").Select()
